# Auto-generated files on 2025-09-07
# Update the HotStock_Top20 worksheet values to reflect the latest rankings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for A2:C21 (row 1 header is unchanged)
$data = @(
    @("先导智能", "先导智能", "卧龙电驱"),
    @("阳光电源", "天赐材料", "先导智能"),
    @("天赐材料", "赣锋锂业", "国轩高科"),
    @("胜宏科技", "露笑科技", "天赐材料"),
    @("赣锋锂业", "阳光电源", "胜宏科技"),
    @("露笑科技", "国轩高科", "吉视传媒"),
    @("天际股份", "亿纬锂能", "阳光电源"),
    @("卧龙电驱", "卧龙电驱", "亿纬锂能"),
    @("国轩高科", "锦浪科技", "露笑科技"),
    @("亿纬锂能", "多氟多", "岩山科技"),
    @("上海电力", "天通股份", "天通股份"),
    @("吉视传媒", "胜宏科技", "赣锋锂业"),
    @("多氟多", "璞泰来", "天际股份"),
    @("天通股份", "天际股份", "金发科技"),
    @("宁德时代", "欣旺达", "上海电力"),
    @("欣旺达", "东方财富", "多氟多"),
    @("锦浪科技", "吉视传媒", "长城军工"),
    @("璞泰来", "当升科技", "中际旭创"),
    @("岩山科技", "上海洗霸", "万通发展"),
    @("上海洗霸", "上海电力", "欣旺达")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
